$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date (Overview) / Correspond Handoff Datetime (de-de)
# -> these share the same text value, update both rows (1c9b3547 + 7e21c464) in both sheets
$wsOverview.Range("G3").Value = "2016-08-21 18:13:56"
$wsOverview.Range("G4").Value = "2016-08-21 18:13:56"
$wsDeDe.Range("H3").Value = "2016-08-21 18:13:56"
$wsDeDe.Range("H4").Value = "2016-08-21 18:13:56"

# zh-cn Priority: ht -> mt (rows 3 and 4)
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# zh-cn Correspond Handoff Datetime (rows 3 and 4)
$wsZhCn.Range("H3").Value = "2016-08-21 18:13:51"
$wsZhCn.Range("H4").Value = "2016-08-21 18:13:51"

# zh-cn Correspond Handback DateTime (rows 3 and 4)
$wsZhCn.Range("K3").Value = "2016-08-21 18:14:13"
$wsZhCn.Range("K4").Value = "2016-08-21 18:14:13"

# de-de Correspond Handback DateTime (rows 3 and 4)
$wsDeDe.Range("K3").Value = "2016-08-21 18:14:19"
$wsDeDe.Range("K4").Value = "2016-08-21 18:14:19"
